# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (interest count) values in column F on the
# "展览" sheet and the mirrored rows on the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$exhibitionUpdates = @{
    2  = 7
    3  = 12861
    4  = 27
    5  = 84
    6  = 68
    7  = 41
    10 = 12780
    11 = 279
    12 = 29
    13 = 8661
    14 = 7662
    15 = 189
    16 = 92
    24 = 16
    25 = 84
}

$allTypesUpdates = @{
    3  = 7
    4  = 12861
    5  = 27
    6  = 84
    7  = 68
    8  = 41
    11 = 12780
    12 = 279
    13 = 29
    14 = 8661
    15 = 7662
    16 = 189
    17 = 92
    26 = 16
    27 = 84
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
